$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.121.28"
$ws.Range("E2").Value = "  -2.65%  "

$ws.Range("D3").Value = "1.847.18"
$ws.Range("E3").Value = "  -1.48%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "0.6935"
$ws.Range("E5").Value = "  -6.32%  "

$ws.Range("D6").Value = "237.93"
$ws.Range("E6").Value = "  -1.93%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.3048"
$ws.Range("E8").Value = "  -3.12%  "

$ws.Range("D9").Value = "0.07602"
$ws.Range("E9").Value = "  +5.32%  "

$ws.Range("D10").Value = "23.39"
$ws.Range("E10").Value = "  -5.02%  "

$ws.Range("D11").Value = "0.08105"
$ws.Range("E11").Value = "  -2.76%  "

$ws.Range("D12").Value = "1.869.03"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("D13").Value = "0.7222"
$ws.Range("E13").Value = "  -3.64%  "

$ws.Range("D14").Value = "5.167"
$ws.Range("E14").Value = "  -4.07%  "

$ws.Range("D15").Value = "88.98"
$ws.Range("E15").Value = "  -3.55%  "

$ws.Range("D16").Value = "29.223.29"
$ws.Range("E16").Value = "  -2.31%  "

$ws.Range("D17").Value = "5.774"
$ws.Range("E17").Value = "  -5.61%  "

$ws.Range("D18").Value = "241.92"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").Value = "0.000007707"
$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").Value = "13.08"
$ws.Range("E20").Value = "  -3.64%  "

$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").Value = "2.108.62"
$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "7.632"
$ws.Range("E24").Value = "  -4.66%  "

$ws.Range("D25").Value = "9.002"
$ws.Range("E25").Value = "  -3.08%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "161.06"
$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "0.1454"
$ws.Range("E27").Value = "  -5.96%  "

$ws.Range("D28").Value = "18.07"
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("D29").Value = "1.934"
$ws.Range("E29").Value = "  -4.26%  "

$ws.Range("D30").Value = "1.388"
$ws.Range("E30").Value = "  -7.59%  "

$ws.Range("D31").Value = "4.414"
$ws.Range("E31").Value = "  -4.14%  "

$ws.Range("E32").Value = "  -2.85%  "

$ws.Range("D33").Value = "4.051"
$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("E34").Value = "  -2.03%  "

$ws.Range("D35").Value = "1.187"
$ws.Range("E35").Value = "  -3.80%  "

$ws.Range("D36").Value = "0.7115"
$ws.Range("E36").Value = "  -5.05%  "

$ws.Range("D37").Value = "1.004"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("E38").Value = "  -1.54%  "

$ws.Range("D39").Value = "0.01859"
$ws.Range("E39").Value = "  -5.11%  "

$ws.Range("D40").Value = "2.691"
$ws.Range("E40").Value = "  -2.22%  "

$ws.Range("D41").Value = "0.9155"
$ws.Range("E41").Value = "  +6.29%  "

$ws.Range("D42").Value = "5.955"
$ws.Range("E42").Value = "  -2.99%  "

$ws.Range("D43").Value = "0.4288"
$ws.Range("E43").Value = "  -5.05%  "

$ws.Range("D44").Value = "69.75"
$ws.Range("E44").Value = "  -3.81%  "

$ws.Range("D45").Value = "1.040.16"
$ws.Range("E45").Value = "  -6.67%  "

$ws.Range("D46").Value = "0.9998"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("D47").Value = "102.60"
$ws.Range("E47").Value = "  -1.62%  "

$ws.Range("D48").Value = "7.213"
$ws.Range("E48").Value = "  -5.23%  "

$ws.Range("D49").Value = "2.014.09"
$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("D50").Value = "1.746"
$ws.Range("E50").Value = "  -6.27%  "

$ws.Range("D51").Value = "9.237"
$ws.Range("E51").Value = "  -2.96%  "
